$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 123
$ws.Range("C5").Value = "dfr"
$ws.Range("D3").Value = "trse"
$ws.Range("B3").Value = "dbs"

$ws.Range("F10").Select()
